# Update column A (index) in the StockEnv 2330/2017 worksheet so that,
# starting from February 2017 (row 18), the index becomes a continuous
# running count instead of restarting at 0 every month.
#
# Rows 2-17  (January 2017) keep their original values (0-15).
# Rows 18-247 get A = (row number - 2), i.e. a running index that
# continues counting across month boundaries for the rest of the year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 18; $r -le 247; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
